# Coursera Capstone Report - "Review Update, Initial Presentation Upload"
#
# 1. Insert a new paragraph ("For a better picture of demographics...") right
#    before the "In reviewing the process" paragraph.
# 2. Split the run in the final ("Going forward...") paragraph that contains
#    "...CMAP only covers the Chicagoland area..." into two runs, with a
#    <w:lastRenderedPageBreak/> on the second one.
# 3. Move the _GoBack bookmark from the end of the document (end of the
#    "Going forward..." paragraph) onto the end of the new paragraph added
#    in step 1.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Step 1: locate the "In reviewing the process" paragraph and push a new,
# blank paragraph in front of it (this shifts it down by one, leaving its
# own XML completely untouched).
# ---------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$reviewParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*In reviewing the process*") {
        $reviewParaIndex = $i
        break
    }
}

$reviewPara = $d.Paragraphs.Item($reviewParaIndex)
$insertPoint = $reviewPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()

# The now-empty paragraph occupies the old index; fill it with the new text
# (the _GoBack bookmark is added separately in step 3, once the original one
# down at the end of the document has been removed).
$newPara = $d.Paragraphs.Item($reviewParaIndex)
$newParaBody = '<w:p><w:r><w:tab/><w:t xml:space="preserve">For a better picture of demographics as they relate to the business, the business could implement a reward system. This reward system could gather age and gender demographics, while also allowing us to see which demographics are most likely to become return customers.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($pkgHeader + $newParaBody + $pkgFooter)

# ---------------------------------------------------------------------
# Step 2: remove the existing _GoBack bookmark (currently sitting at the
# end of the document) before touching the last paragraph, so the XML
# rewrite below does not resurrect it around the new runs.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Step 3: rewrite the final ("Going forward...") paragraph, splitting the
# run that contains "...CMAP only covers the Chicagoland area..." so a
# <w:lastRenderedPageBreak/> sits between the two halves. Exclude the very
# last character of the range (the document's closing paragraph mark) so
# Word edits this paragraph in place instead of appending a fresh empty
# one.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRng = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$goingForwardBody = '<w:p><w:r><w:tab/><w:t xml:space="preserve">Going forward, this methodology </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>could be applied</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to a program to allow any business owner to determine where to expand. Doing so would require changing the demographic data source, as CMAP </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>only covers the Chicagoland area, and not the United Stat</w:t></w:r><w:r><w:t xml:space="preserve">es as a whole. The U.S. Census might be a viable option for such an application. </w:t></w:r></w:p>'
$lastRng.InsertXML($pkgHeader + $goingForwardBody + $pkgFooter)

# ---------------------------------------------------------------------
# Step 4: re-create the _GoBack bookmark at the end of the new paragraph
# inserted in step 1 (collapsed, right after its text, same as the
# original placement at the end of the document).
# ---------------------------------------------------------------------
$newPara = $d.Paragraphs.Item($reviewParaIndex)
$bookmarkPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null

Write-Output "done"
